# fix: alterar python version para 3.11.5
# Update absenteeism_data records (rows 2-11) with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 13604
$ws.Range("B2").Value = "Ana Laura Rezende"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45081
$ws.Range("G2").Value = 5337.07

# Row 3
$ws.Range("A3").Value = 33493
$ws.Range("B3").Value = "Dra. Emilly Viana"
$ws.Range("C3").Value = "Operações"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 45101
$ws.Range("G3").Value = 11403.17

# Row 4
$ws.Range("A4").Value = 28519
$ws.Range("B4").Value = "Dr. Carlos Eduardo Almeida"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 45095
$ws.Range("G4").Value = 8914.16

# Row 5
$ws.Range("A5").Value = 14164
$ws.Range("B5").Value = "Srta. Lívia Ferreira"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45105
$ws.Range("G5").Value = 11422.44

# Row 6
$ws.Range("A6").Value = 13163
$ws.Range("B6").Value = "Mariana Costa"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45087
$ws.Range("G6").Value = 11826.01

# Row 7
$ws.Range("A7").Value = 23224
$ws.Range("B7").Value = "Igor Araújo"
$ws.Range("C7").Value = "Jurídico"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45079
$ws.Range("G7").Value = 5419.78

# Row 8
$ws.Range("A8").Value = 58311
$ws.Range("B8").Value = "Davi Luiz Nascimento"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Doença"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 45085
$ws.Range("G8").Value = 11545.2

# Row 9
$ws.Range("A9").Value = 54724
$ws.Range("B9").Value = "Enzo Gabriel Martins"
$ws.Range("C9").Value = "Engenharia"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 45105
$ws.Range("G9").Value = 3483.1

# Row 10
$ws.Range("A10").Value = 24559
$ws.Range("B10").Value = "João Felipe Azevedo"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 45101
$ws.Range("G10").Value = 8919.29

# Row 11
$ws.Range("A11").Value = 38730
$ws.Range("B11").Value = "João Pedro da Rosa"
$ws.Range("C11").Value = "Vendas"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45089
$ws.Range("G11").Value = 2707.2
